$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "54.732.20"
Set-TextValue "E2" "  +7.33%  "
Set-TextValue "D3" "2.428.09"
Set-TextValue "E3" "  +8.12%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "478.66"
Set-TextValue "E5" "  +11.93%  "
Set-TextValue "D6" "138.66"
Set-TextValue "E6" "  +21.79%  "
Set-TextValue "D7" "0.996"
Set-TextValue "E7" "  -0.25%  "
Set-TextValue "E8" "  +11.65%  "
Set-TextValue "D9" "2.457.22"
Set-TextValue "E9" "  +9.22%  "
Set-TextValue "D10" "0.0956"
Set-TextValue "E10" "  +16.39%  "
Set-TextValue "E11" "  +6.94%  "
Set-TextValue "E12" "  +10.59%  "
Set-TextValue "E13" "  +2.51%  "
Set-TextValue "D14" "2.868.67"
Set-TextValue "E14" "  +8.77%  "
Set-TextValue "D15" "54.874.12"
Set-TextValue "E15" "  +7.45%  "
Set-TextValue "D16" "20.44"
Set-TextValue "E16" "  +13.28%  "
Set-TextValue "D17" "0.0000133"
Set-TextValue "E17" "  +20.77%  "
Set-TextValue "D18" "2.454.24"
Set-TextValue "E18" "  +8.00%  "
Set-TextValue "E19" "  +13.51%  "
Set-TextValue "D20" "9.92"
Set-TextValue "E20" "  +18.83%  "
Set-TextValue "D21" "312.28"
Set-TextValue "E21" "  +8.41%  "
Set-TextValue "D22" "0.996"
Set-TextValue "E22" "  +0.23%  "
Set-TextValue "D23" "5.64"
Set-TextValue "E23" "  +15.68%  "
Set-TextValue "D24" "57.06"
Set-TextValue "E24" "  +9.36%  "
Set-TextValue "E25" "  +1.64%  "
Set-TextValue "E26" "  +13.37%  "
Set-TextValue "E27" "  +18.82%  "
Set-TextValue "D28" "2.559.67"
Set-TextValue "E28" "  +9.41%  "
Set-TextValue "D29" "7.32"
Set-TextValue "E29" "  +11.41%  "
Set-TextValue "D30" "0.0₃0767"
Set-TextValue "E30" "  +26.81%  "
Set-TextValue "D31" "0.998"
Set-TextValue "E31" "  +0.00%  "
Set-TextValue "D32" "148.41"
Set-TextValue "E32" "  +4.30%  "
Set-TextValue "D33" "17.86"
Set-TextValue "E33" "  +10.61%  "
Set-TextValue "E34" "  +15.35%  "
Set-TextValue "D35" "5.15"
Set-TextValue "E35" "  +14.12%  "
Set-TextValue "D36" "1.12"
Set-TextValue "E36" "  +18.11%  "
Set-TextValue "E37" "  +11.35%  "
Set-TextValue "D38" "3.58"
Set-TextValue "E38" "  +10.34%  "
Set-TextValue "B39" "OKB"
Set-TextValue "C39" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D39" "33.36"
Set-TextValue "E39" "  +6.21%  "
Set-TextValue "B40" "FirstDigitalUSD"
Set-TextValue "C40" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D40" "0.991"
Set-TextValue "E40" "  -0.64%  "
Set-TextValue "D41" "0.602"
Set-TextValue "E41" "  +10.84%  "
Set-TextValue "D42" "3.40"
Set-TextValue "E42" "  +13.63%  "
Set-TextValue "D43" "0.0540"
Set-TextValue "E43" "  +12.16%  "
Set-TextValue "D44" "1.29"
Set-TextValue "E44" "  +16.32%  "
Set-TextValue "D45" "10.11"
Set-TextValue "E45" "  -0.02%  "
Set-TextValue "D46" "255.41"
Set-TextValue "E46" "  +35.43%  "
Set-TextValue "D47" "4.62"
Set-TextValue "E47" "  +21.69%  "
Set-TextValue "D48" "0.0889"
Set-TextValue "E48" "  +13.31%  "
Set-TextValue "D49" "1.922.81"
Set-TextValue "E49" "  +4.47%  "
Set-TextValue "D50" "0.0221"
Set-TextValue "E50" "  +12.96%  "
Set-TextValue "D51" "17.09"
Set-TextValue "E51" "  +14.10%  "
